# Dedicate only one column for a dependee in the dependency section.
#
# The sheet lays out a dependency tree where each nesting level used to get
# its own column (C/D for the arrow markers, D/E/F for the labels). This
# collapsed one redundant indent level: column D is removed entirely, so
# everything that lived in column E shifts into D, and everything that
# lived in column F shifts into D's old neighbour E.
#
# One row (18) had its left-pointing arrow sitting in the column being
# removed (D18) with nothing in the column to its left (C18 was blank) -
# deleting the column would otherwise silently drop that arrow, so it is
# re-added immediately to the left of the label it now points at (C18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column D; Excel shifts everything right of it (E, F, ...) one
# column to the left, taking each cell's value/style along with it.
$ws.Columns("D").Delete()

# Row 18's "<-" marker lived in the deleted column (D18) and has no
# survivor to its left, so restore it one column left of the label (D18,
# formerly E18) it points to.
$ws.Range("C18").Value = "←"
